$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.400.19'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.817.61'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.31'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5114'
$ws.Range('E7').Value = '  -4.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3947'
$ws.Range('E8').Value = '  -2.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08051'
$ws.Range('E9').Value = '  +5.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.68'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('E12').Value = '  +0.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.260'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.497'
$ws.Range('E15').Value = '  -1.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.816.08'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('E17').Value = '  +6.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.65'
$ws.Range('E18').Value = '  +3.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06634'
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.67'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.086'
$ws.Range('E22').Value = '  -0.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.436.92'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.26'
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.269'
$ws.Range('E25').Value = '  +3.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.13'
$ws.Range('E26').Value = '  +2.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.033.11'
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '154.94'
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.402'
$ws.Range('E29').Value = '  -2.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.90'
$ws.Range('E30').Value = '  +1.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1101'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.775'
$ws.Range('E33').Value = '  +2.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.652'
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07023'
$ws.Range('E35').Value = '  -3.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2224'
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02325'
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.207'
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.811'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6258'
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.397'
$ws.Range('E44').Value = '  +0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.50'
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5907'
$ws.Range('E47').Value = '  +1.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.89'
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('E50').Value = '  -1.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06887'
$ws.Range('E51').Value = '  +0.05%  '
